$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 94

$ws.Cells.Item($row, 1).Value = 45457.2916666667
$ws.Cells.Item($row, 2).Value = 500
$ws.Cells.Item($row, 3).Value = 0.709999978542328
$ws.Cells.Item($row, 4).Value = 0.704999983310699
$ws.Cells.Item($row, 5).Value = 0.709999978542328
$ws.Cells.Item($row, 6).Value = 0.704999983310699
$ws.Cells.Item($row, 7).NumberFormat = "@"
$ws.Cells.Item($row, 7).Value = "0.704999983310699"
$ws.Cells.Item($row, 7).Style = "Normal"
$ws.Cells.Item($row, 8).Value = "BWZ.MI"

$ws.Cells.Item($row, 1).Style = $ws.Cells.Item($row - 1, 1).Style
